$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header / table column from "氏名" to "教員名"
$ws.Range("A1").Value = "教員名"

# Move the selection/active cell to A1 (matches the saved view state)
$ws.Range("A1").Select()
